$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Formula = "'247.47"
$ws.Range("D2").Style = "Normal"

# Row 3
$ws.Range("D3").Formula = "'22.39"
$ws.Range("D3").Style = "Normal"

# Row 4
$ws.Range("D4").Formula = "'5.523"
$ws.Range("D4").Style = "Normal"

# Row 5
$ws.Range("D5").Formula = "'0.05621"
$ws.Range("D5").Style = "Normal"

# Row 6
$ws.Range("B6").Value = 'KuCoinToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D6").Formula = "'6.474"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '5KuCoinTokenKCS'

# Row 7
$ws.Range("B7").Value = 'MXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D7").Formula = "'0.8048"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '6MXTokenMX'

# Row 8
$ws.Range("B8").Value = 'FTXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D8").Formula = "'1.048"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '7FTXTokenFTT'

# Row 9
$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9").Formula = "'0.1423"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '8WazirXWRX'

# Row 10
$ws.Range("B10").Value = 'MandalaExchangeToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D10").Formula = "'0.07272"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '9MandalaExchangeTokenMDX'

# Row 11
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").Formula = "'0.03190"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '10LiechtensteinCryptoassetsExchangeLCX'

# Row 12
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Formula = "'0.02959"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '11BitrueCoinBTR'

# Row 13
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Formula = "'0.09266"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '12BitMartTokenBMX'

# Row 14
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Formula = "'0.001661"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '13BitForexTokenBF'

# Row 15
$ws.Range("B15").Value = 'MCDex'
$ws.Range("C15").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D15").Formula = "'3.199"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '14MCDexMCB'

# Row 16
$ws.Range("B16").Value = 'CoinExToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D16").Formula = "'0.04699"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '15CoinExTokenCET'

# Row 17
$ws.Range("B17").Value = 'One'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D17").Formula = "'0.0005963"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '16OneONE'

# Row 18
$ws.Range("B18").Value = 'TigerCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D18").Formula = "'0.006287"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '17TigerCashTCH'

# Row 19
$ws.Range("B19").Value = 'BitKan'
$ws.Range("C19").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D19").Formula = "'0.001049"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '18BitKanKAN'

# Row 20
$ws.Range("B20").Value = 'HotbitToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D20").Formula = "'0.003816"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '19HotbitTokenHTB'

# Row 21
$ws.Range("B21").Value = 'NitroEx'
$ws.Range("C21").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D21").Formula = "'0.0001502"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '20NitroExNTX'

# Row 22
$ws.Range("B22").Value = 'UpBots'
$ws.Range("C22").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D22").Formula = "'0.0003307"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '21UpBotsUBXT'

# Row 23
$ws.Range("B23").Value = 'LEO'
$ws.Range("C23").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D23").Formula = "'3.972"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '22LEOLEO'

# Row 24
$ws.Range("B24").Value = 'GateToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D24").Formula = "'3.386"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '23GateTokenGT'

# Row 40
$ws.Range("D40").Formula = "'0.04167"
$ws.Range("D40").Style = "Normal"

# Row 41
$ws.Range("D41").Formula = "'0.1043"
$ws.Range("D41").Style = "Normal"

# Row 42
$ws.Range("D42").Formula = "'0.002975"
$ws.Range("D42").Style = "Normal"

# Row 43
$ws.Range("D43").Formula = "'0.006861"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '42KickTokenKICK'

# Row 44
$ws.Range("D44").Formula = "'0.009079"
$ws.Range("D44").Style = "Normal"

# Row 45
$ws.Range("D45").Formula = "'0.00005636"
$ws.Range("D45").Style = "Normal"

# Row 46
$ws.Range("D46").Formula = "'0.00000000752"
$ws.Range("D46").Style = "Normal"

# Row 47
$ws.Range("D47").Formula = "'0.6815"
$ws.Range("D47").Style = "Normal"

# Row 48
$ws.Range("D48").Formula = "'0.02526"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '47BOLOBOLOWorstin24h'

# Row 49
$ws.Range("D49").Formula = "'0.00002105"
$ws.Range("D49").Style = "Normal"
